$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (currently sitting right after
#    the "third-party" run). Word keeps only one "_GoBack" bookmark,
#    tracking the most recent edit location, so it will be re-created
#    at the new edit location below.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) Rebuild the "Feral Kingdom ..." paragraph:
#      - split "Pokemon-esque" so "esque" sits in its own run,
#        bracketed by spellcheck proofErr markers (as Word does when
#        a word is retyped and flagged by the spell checker)
#      - append the new sentence about ending a battle
#      - drop a fresh "_GoBack" bookmark at the end of the new text
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Feral Kingdom will be a")
$para = $target.Paragraphs(1).Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1486F145" w14:textId="205FCDAF" w:rsidR="00DF4E0D" w:rsidRPr="00DF4E0D" w:rsidRDefault="003858F1" w:rsidP="00DF4E0D"><w:pPr><w:ind w:left="360"/><w:rPr><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>Feral Kingdom will be a Pokémon-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>esque</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> turn based combat game with </w:t></w:r><w:r w:rsidR="003709DE"><w:rPr><w:sz w:val="20"/></w:rPr><w:t>an</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> overworld to move between battles</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">, the player will  be able to complete a battle either by killing the opposing monster in the battle or by running away from the battle, the latter </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$para.InsertXML($xml)

# ------------------------------------------------------------------
# 3) Bump the DATE field result shown in the header.
# ------------------------------------------------------------------
$header = $d.Sections(1).Headers(1)
$header.Range.Find.Execute("19/05/2020", $true, $false, $false, $false, $false,
                            $true, 1, $false, "22/05/2020", 2)
